# Auto-generated COM-interop script: adds a new 'Ready for handoff' row
# for file 96092630-0683-... across Overview / zh-cn / de-de sheets,
# mirroring the existing 9b63b6e9-... row that is already present.

$wb = $excel.ActiveWorkbook

$fileNameNew = '96092630-0683-43d1-9e3f-5a526ba3fb5fooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$pathAndNameNew = 'e2e\96092630-0683-43d1-9e3f-5a526ba3fb5fooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$statusNew = 'Ready for handoff'
$dateOverview = '2016-08-23 16:29:42'
$zhCnXlfNew = '96092630-0683-43d1-9e3f-5a526ba3fb5foooooooooooooooooooooooooooooooooooooooo.fe50d7c0e5d2a8389431a7fb99fe5452becd0c16.zh-cn.xlf'
$dateZh = '2016-08-23 16:29:37'
$deDeXlfNew = '96092630-0683-43d1-9e3f-5a526ba3fb5foooooooooooooooooooooooooooooooooooooooo.fe50d7c0e5d2a8389431a7fb99fe5452becd0c16.de-de.xlf'
$hyperlinkTarget = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcab026878ef5e17da3d8771f951f4031af3f72d/e2e/96092630-0683-43d1-9e3f-5a526ba3fb5fooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

# ---- Overview sheet: append row referencing the new file ----
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A3").Value = $fileNameNew
$wsOverview.Range("B3").Value = $pathAndNameNew
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $dateOverview
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $pathAndNameNew) | Out-Null

# ---- zh-cn sheet: append row referencing the new file ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null
$wsZhCn.Range("A3").Value = $fileNameNew
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhCnXlfNew
$wsZhCn.Range("H3").Value = $dateZh
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $fileNameNew) | Out-Null

# ---- de-de sheet: append row referencing the new file ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null
$wsDeDe.Range("A3").Value = $fileNameNew
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $deDeXlfNew
$wsDeDe.Range("H3").Value = $dateOverview
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $fileNameNew) | Out-Null

